$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2101.617
$ws.Range("J17").Value = 2101.617
$ws.Range("L17").Value = 6304.851000000001
$ws.Range("N17").Value = -6640.851000000001

$ws.Range("H28").Value = 2958.3
$ws.Range("J28").Value = 1998
$ws.Range("L28").Value = 1998
$ws.Range("N28").Value = -2968

$ws.Range("H58").Value = 7540.143
$ws.Range("I58").Value = 312.6
$ws.Range("J58").Value = 11555.444
$ws.Range("K58").Value = 937.8000000000001
$ws.Range("L58").Value = 34666.33199999999
$ws.Range("M58").Value = -787.8000000000001
$ws.Range("N58").Value = -34966.33199999999

$ws.Range("H98").Value = 1680.3077
$ws.Range("I98").Value = 1653.6945
$ws.Range("K98").Value = 1653.6945
$ws.Range("M98").Value = -155.6945000000001

$ws.Range("H112").Value = 1533.0256
$ws.Range("I112").Value = 971.5
$ws.Range("J112").Value = 1563.3784
$ws.Range("K112").Value = 2914.5
$ws.Range("L112").Value = 4690.135200000001
$ws.Range("M112").Value = -1806.5
$ws.Range("N112").Value = -6906.135200000001

$ws.Range("H116").Value = 2315.8076
$ws.Range("I116").Value = 2378.5715
$ws.Range("J116").Value = 2242.5833
$ws.Range("K116").Value = 2378.5715
$ws.Range("L116").Value = 2242.5833
$ws.Range("M116").Value = 1063.4285
$ws.Range("N116").Value = -9126.5833

$ws.Range("H122").Value = 1680.3077
$ws.Range("I122").Value = 1653.6945
$ws.Range("K122").Value = 4961.083500000001
$ws.Range("M122").Value = -2511.083500000001

$ws.Range("H138").Value = 18184162
$ws.Range("J138").Value = 3118.0645
$ws.Range("L138").Value = 9354.193499999999
$ws.Range("N138").Value = -19634.1935

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1751.0461
$ws.Range("I32").Value = 1664.6031
$ws.Range("K32").Value = 1664.6031
$ws.Range("M32").Value = -1377.6031

$ws.Range("H74").Value = 3344.4897
$ws.Range("I74").Value = 2893.2942
$ws.Range("J74").Value = 4367.2
$ws.Range("K74").Value = 2893.2942
$ws.Range("L74").Value = 4367.2
$ws.Range("M74").Value = -2019.2942
$ws.Range("N74").Value = -6115.2

$ws.Range("H77").Value = 3344.4897
$ws.Range("I77").Value = 2893.2942
$ws.Range("J77").Value = 4367.2
$ws.Range("K77").Value = 14466.471
$ws.Range("L77").Value = 21836
$ws.Range("M77").Value = -10098.471
$ws.Range("N77").Value = -30572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3834.4736
$ws.Range("I86").Value = 3185.4443
$ws.Range("J86").Value = 4418.6
$ws.Range("K86").Value = 3185.4443
$ws.Range("L86").Value = 4418.6
$ws.Range("M86").Value = -2062.4443
$ws.Range("N86").Value = -6664.6

$ws.Range("H89").Value = 3834.4736
$ws.Range("I89").Value = 3185.4443
$ws.Range("J89").Value = 4418.6
$ws.Range("K89").Value = 15927.2215
$ws.Range("L89").Value = 22093
$ws.Range("M89").Value = -10311.2215
$ws.Range("N89").Value = -33325

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 6255
$ws.Range("I4").Value = 10
$ws.Range("K4").Value = 10
$ws.Range("M4").Value = 102

$ws.Range("H99").Value = 4994.222
$ws.Range("I99").Value = 2901.3635
$ws.Range("J99").Value = 8283
$ws.Range("K99").Value = 2901.3635
$ws.Range("L99").Value = 8283
$ws.Range("M99").Value = -1403.3635
$ws.Range("N99").Value = -11279

$ws.Range("H126").Value = 4994.222
$ws.Range("I126").Value = 2901.3635
$ws.Range("J126").Value = 8283
$ws.Range("K126").Value = 8704.0905
$ws.Range("L126").Value = 24849
$ws.Range("M126").Value = -6234.0905
$ws.Range("N126").Value = -29789

$ws.Range("H134").Value = 1980.8422
$ws.Range("I134").Value = 1952.5294
$ws.Range("K134").Value = 5857.5882
$ws.Range("M134").Value = -3322.5882

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 75046.71000000001
$ws.Range("I4").Value = 80127.08
$ws.Range("K4").Value = 240381.24
$ws.Range("M4").Value = -240269.24

$ws.Range("H98").Value = 865.3333
$ws.Range("I98").Value = 750
$ws.Range("J98").Value = 923
$ws.Range("K98").Value = 2250
$ws.Range("L98").Value = 2769
$ws.Range("M98").Value = -752
$ws.Range("N98").Value = -5765

$ws.Range("H131").Value = 1410.8
$ws.Range("J131").Value = 1818
$ws.Range("L131").Value = 5454
$ws.Range("N131").Value = -15534

$ws.Range("H133").Value = 4710.5884
$ws.Range("I133").Value = 4007.2727
$ws.Range("K133").Value = 12021.8181
$ws.Range("M133").Value = -6961.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H11").Value = 423996.25
$ws.Range("J11").Value = 540328.3
$ws.Range("L11").Value = 540328.3
$ws.Range("N11").Value = -540606.3

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H113").Value = 1466.9333
$ws.Range("I113").Value = 626
$ws.Range("K113").Value = 626
$ws.Range("M113").Value = 1544

$ws.Range("H122").Value = 2635.84
$ws.Range("J122").Value = 3781.6667
$ws.Range("L122").Value = 11345.0001
$ws.Range("N122").Value = -16245.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6145.846
$ws.Range("I7").Value = 4234.875
$ws.Range("K7").Value = 4234.875
$ws.Range("M7").Value = -4122.875

$ws.Range("H16").Value = 759.5
$ws.Range("I16").Value = 806.3182
$ws.Range("K16").Value = 806.3182
$ws.Range("M16").Value = -636.3182

$ws.Range("H40").Value = 3742.262
$ws.Range("I40").Value = 3700.282
$ws.Range("J40").Value = 4288
$ws.Range("K40").Value = 3700.282
$ws.Range("L40").Value = 4288
$ws.Range("M40").Value = -3564.282
$ws.Range("N40").Value = -4560

$ws.Range("H61").Value = 1903.1818
$ws.Range("I61").Value = 1562.2858
$ws.Range("J61").Value = 2499.75
$ws.Range("K61").Value = 1562.2858
$ws.Range("L61").Value = 2499.75
$ws.Range("M61").Value = -1360.2858
$ws.Range("N61").Value = -2903.75

$ws.Range("H100").Value = 5342.125
$ws.Range("I100").Value = 3236
$ws.Range("J100").Value = 9975.6
$ws.Range("K100").Value = 3236
$ws.Range("L100").Value = 9975.6
$ws.Range("M100").Value = -2695
$ws.Range("N100").Value = -11057.6

$ws.Range("H113").Value = 1903.1818
$ws.Range("I113").Value = 1562.2858
$ws.Range("J113").Value = 2499.75
$ws.Range("K113").Value = 1562.2858
$ws.Range("L113").Value = 2499.75
$ws.Range("M113").Value = 607.7141999999999
$ws.Range("N113").Value = -6839.75

$ws.Range("H126").Value = 6145.846
$ws.Range("I126").Value = 4234.875
$ws.Range("K126").Value = 12704.625
$ws.Range("M126").Value = -10234.625

$ws.Range("H132").Value = 3243.36
$ws.Range("I132").Value = 3349.8472
$ws.Range("J132").Value = 2639.9333
$ws.Range("K132").Value = 10049.5416
$ws.Range("L132").Value = 7919.7999
$ws.Range("M132").Value = -7519.5416
$ws.Range("N132").Value = -12979.7999

$ws.Range("H136").Value = 5739.9375
$ws.Range("I136").Value = 4995.3335
$ws.Range("K136").Value = 14986.0005
$ws.Range("M136").Value = -12436.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 20000
$ws.Range("I2").Value = 20000
$ws.Range("K2").Value = 20000
$ws.Range("M2").Value = -19888

$ws.Range("H41").Value = 20540.715
$ws.Range("I41").Value = 21051.666
$ws.Range("J41").Value = 20157.5
$ws.Range("K41").Value = 21051.666
$ws.Range("L41").Value = 20157.5
$ws.Range("M41").Value = -20661.666
$ws.Range("N41").Value = -20937.5

$ws.Range("H104").Value = 19456.666
$ws.Range("J104").Value = 19456.666
$ws.Range("L104").Value = 19456.666
$ws.Range("N104").Value = -26444.666
